$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the price/volume columns to Text format before writing values so that
# numeric-looking strings (e.g. "318.89", "0.1684") are not auto-converted to numbers
# by Excel and instead keep their exact original text content (including padding spaces).
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '27.640.33'
$ws.Range("E2").Value = '  -0.50%  '
$ws.Range("D3").Value = '1.842.53'
$ws.Range("E3").Value = '  -1.25%  '
$ws.Range("E4").Value = '  -1.96%  '
$ws.Range("D5").Value = '318.89'
$ws.Range("E5").Value = '  -1.50%  '
$ws.Range("E6").Value = '  -2.10%  '
$ws.Range("E7").Value = '  -2.75%  '
$ws.Range("D8").Value = '0.3739'
$ws.Range("E8").Value = '  -1.98%  '
$ws.Range("D9").Value = '0.07335'
$ws.Range("E9").Value = '  -1.87%  '
$ws.Range("D10").Value = '0.8749'
$ws.Range("E10").Value = '  -1.58%  '
$ws.Range("D11").Value = '21.58'
$ws.Range("E11").Value = '  -0.93%  '
$ws.Range("D12").Value = '1.848.74'
$ws.Range("E12").Value = '  -1.56%  '
$ws.Range("D13").Value = '6.718'
$ws.Range("E13").Value = '  -1.06%  '
$ws.Range("D14").Value = '5.438'
$ws.Range("E14").Value = '  -2.14%  '
$ws.Range("D15").Value = '0.07108'
$ws.Range("E15").Value = '  -1.34%  '
$ws.Range("D16").Value = '89.12'
$ws.Range("E16").Value = '  +5.45%  '
$ws.Range("D17").Value = '1.013'
$ws.Range("E17").Value = '  -2.45%  '
$ws.Range("D18").Value = '0.000008972'
$ws.Range("E18").Value = '  -1.89%  '
$ws.Range("E19").Value = '  -2.10%  '
$ws.Range("D20").Value = '15.45'
$ws.Range("E20").Value = '  -0.95%  '
$ws.Range("D21").Value = '27.650.24'
$ws.Range("E21").Value = '  -0.50%  '
$ws.Range("E22").Value = '  -2.19%  '
$ws.Range("D23").Value = '11.08'
$ws.Range("E23").Value = '  -2.14%  '
$ws.Range("D24").Value = '2.073.92'
$ws.Range("E24").Value = '  -0.86%  '
$ws.Range("D25").Value = '2.005'
$ws.Range("E25").Value = '  -0.75%  '
$ws.Range("D26").Value = '155.66'
$ws.Range("E26").Value = '  -1.88%  '
$ws.Range("E27").Value = '  -1.65%  '
$ws.Range("D28").Value = '2.165'
$ws.Range("E28").Value = '  +8.83%  '
$ws.Range("D29").Value = '5.358'
$ws.Range("E29").Value = '  -0.31%  '
$ws.Range("D30").Value = '118.67'
$ws.Range("E30").Value = '  -0.33%  '
$ws.Range("D31").Value = '0.08938'
$ws.Range("E31").Value = '  -1.34%  '
$ws.Range("D32").Value = '1.227'
$ws.Range("E32").Value = '  -0.70%  '
$ws.Range("D33").Value = '0.7763'
$ws.Range("E33").Value = '  -0.70%  '
$ws.Range("D34").Value = '4.540'
$ws.Range("E34").Value = '  -1.33%  '
$ws.Range("D35").Value = '2.879'
$ws.Range("E35").Value = '  -4.92%  '
$ws.Range("E36").Value = '  -2.17%  '
$ws.Range("D37").Value = '1.135'
$ws.Range("E37").Value = '  -0.97%  '
$ws.Range("D38").Value = '0.05328'
$ws.Range("E38").Value = '  -0.78%  '
$ws.Range("D39").Value = '0.01973'
$ws.Range("E39").Value = '  -1.04%  '
$ws.Range("D40").Value = '7.294'
$ws.Range("E40").Value = '  +5.26%  '
$ws.Range("D41").Value = '2.935'
$ws.Range("E41").Value = '  +1.86%  '
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").Value = '0.1684'
$ws.Range("E42").Value = '  -0.68%  '
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").Value = '0.5108'
$ws.Range("E43").Value = '  -2.27%  '
$ws.Range("E44").Value = '  +0.39%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '10.73'
$ws.Range("E45").Value = '  +0.16%  '
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").Value = '109.11'
$ws.Range("E46").Value = '  -2.35%  '
$ws.Range("D47").Value = '0.4740'
$ws.Range("E47").Value = '  +0.03%  '
$ws.Range("D48").Value = '0.06476'
$ws.Range("E48").Value = '  -2.97%  '
$ws.Range("E49").Value = '  -2.33%  '
$ws.Range("D50").Value = '1.689'
$ws.Range("E50").Value = '  -1.61%  '
$ws.Range("D51").Value = '1.839'
$ws.Range("E51").Value = '  -4.59%  '

# Restore the original (Normal) style on the range so the cells keep the same
# appearance/format as before the edit (no leftover Text number format).
$dataRange.Style = "Normal"

